# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") was regenerated from the live
# source data (strikeout counts replacing the prior "Strike#" derived
# values). This updates each data row (2-76) in column G to the newly
# computed K value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 3
    16 = 0
    17 = 0
    18 = 2
    19 = 1
    20 = 2
    21 = 3
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 1
    50 = 2
    51 = 3
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 1
    61 = 4
    62 = 0
    63 = 1
    64 = 5
    65 = 1
    66 = 0
    67 = 0
    68 = 0
    69 = 2
    70 = 2
    71 = 0
    72 = 1
    73 = 2
    74 = 1
    75 = 1
    76 = 2
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
